# Weekly fruit/vegetable price update: insert a new weekly record for
# "Apio" (Terminal La Palmera de La Serena) as row 645, pushing the
# existing rows (old 645-696) down by one (new 646-697).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 645 (shifts rows 645..696 -> 646..697,
# carrying their values/styles with them, matching the target diff).
$ws.Rows.Item(645).Insert()

# Populate the newly inserted row 645 with the new weekly data point.
$ws.Cells.Item(645, 1).Value = 8
$ws.Cells.Item(645, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(645, 3).Value = "Coquimbo"
$ws.Cells.Item(645, 4).Value = 45166
$ws.Cells.Item(645, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(645, 5).Value = 4
$ws.Cells.Item(645, 6).Value = 100112017
$ws.Cells.Item(645, 7).Value = "Apio"
$ws.Cells.Item(645, 8).Value = "Americana (o)"
$ws.Cells.Item(645, 9).Value = "Primera"
$ws.Cells.Item(645, 10).Value = 1000
$ws.Cells.Item(645, 11).Value = 7000
$ws.Cells.Item(645, 12).Value = 8000
$ws.Cells.Item(645, 13).Value = 7500
$ws.Cells.Item(645, 14).Value = "$/docena de matas"
$ws.Cells.Item(645, 15).Value = "Provincia del Elqu"+[char]0x00ED
$ws.Cells.Item(645, 16).Value = 1250
$ws.Cells.Item(645, 17).Value = 6
$ws.Cells.Item(645, 18).Value = "Hortaliza"
